# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
$headerRange = $ws.Range("AD1:AF1")

# Match the formatting used by the other header cells (e.g. AC1):
# bold font, thin border all around, centered/top aligned.
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# --- Data rows (2-56): every team/player row gets the same record ---
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 80   # column AD
    $ws.Cells.Item($r, 31).Value2 = 82   # column AE
    $ws.Cells.Item($r, 32).Value2 = 0    # column AF
}
